$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column C
$ws.Range("C1").Value = "z"

# Row 2 is a plain (non-shared) formula
$ws.Range("C2").Formula = "=A2^3-5*A2^2+500000"

# Rows 3-10 use a shared formula (C3 is the "source" with si="1")
$ws.Range("C3:C10").Formula = "=A3^3-5*A3^2+500000"

# Set column C width to fit contents (bestFit / autofit); stored width rounds to 10
$ws.Columns.Item(3).ColumnWidth = 9.1

# Update the selected cell, as shown in the diff (activeCell moved from E8 to D6)
$ws.Range("D6").Select()

$wb.Save()
